$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "język" (language) column (column F) entirely - shifting
# subsequent columns (grupa, nr tel, email, notatka rekrutacyjna) left.
$ws.Range("F1:F3").EntireColumn.Delete()

# Reset selection to F1 to match the post-edit state
$ws.Range("F1").Select()
